$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1203.2727
$ws.Range("I15").Value = 1203.2727
$ws.Range("K15").Value = 3609.8181
$ws.Range("M15").Value = -3440.8181
$ws.Range("H31").Value = 2490
$ws.Range("I31").Value = 1112.5
$ws.Range("K31").Value = 3337.5
$ws.Range("M31").Value = -3107.5
$ws.Range("H62").Value = 10905.857
$ws.Range("I62").Value = 8395
$ws.Range("J62").Value = 11324.333
$ws.Range("K62").Value = 8395
$ws.Range("L62").Value = 11324.333
$ws.Range("M62").Value = -7771
$ws.Range("N62").Value = -12572.333
$ws.Range("H65").Value = 10905.857
$ws.Range("I65").Value = 8395
$ws.Range("J65").Value = 11324.333
$ws.Range("K65").Value = 41975
$ws.Range("L65").Value = 56621.665
$ws.Range("M65").Value = -38855
$ws.Range("N65").Value = -62861.665
$ws.Range("H138").Value = 6865.523
$ws.Range("I138").Value = 5600
$ws.Range("J138").Value = 6948.5083
$ws.Range("K138").Value = 16800
$ws.Range("L138").Value = 20845.5249
$ws.Range("M138").Value = -11660
$ws.Range("N138").Value = -31125.5249
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 74188.14
$ws.Range("I74").Value = 127385.625
$ws.Range("K74").Value = 127385.625
$ws.Range("M74").Value = -126511.625
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H77").Value = 74188.14
$ws.Range("I77").Value = 127385.625
$ws.Range("K77").Value = 636928.125
$ws.Range("M77").Value = -632560.125
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H122").Value = 2093.625
$ws.Range("I122").Value = 2030.8096
$ws.Range("J122").Value = 2533.3333
$ws.Range("K122").Value = 6092.4288
$ws.Range("L122").Value = 7599.999899999999
$ws.Range("M122").Value = -3642.4288
$ws.Range("N122").Value = -12499.9999
$ws.Range("H133").Value = 177499.5
$ws.Range("J133").Value = 177499.5
$ws.Range("L133").Value = 177499.5
$ws.Range("N133").Value = -182559.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2703.2856
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 4833.9375
$ws.Range("I31").Value = 2144.1428
$ws.Range("J31").Value = 6926
$ws.Range("K31").Value = 2144.1428
$ws.Range("L31").Value = 6926
$ws.Range("M31").Value = -1849.1428
$ws.Range("N31").Value = -7516
$ws.Range("H34").Value = 4833.9375
$ws.Range("I34").Value = 2144.1428
$ws.Range("J34").Value = 6926
$ws.Range("K34").Value = 2144.1428
$ws.Range("L34").Value = 6926
$ws.Range("M34").Value = -1942.1428
$ws.Range("N34").Value = -7330
$ws.Range("H70").Value = 110000
$ws.Range("J70").Value = 110000
$ws.Range("L70").Value = 110000
$ws.Range("N70").Value = -110630
$ws.Range("H73").Value = 110000
$ws.Range("J73").Value = 110000
$ws.Range("L73").Value = 110000
$ws.Range("N73").Value = -112184
$ws.Range("H103").Value = 8005.5
$ws.Range("I103").Value = 8005.5
$ws.Range("K103").Value = 8005.5
$ws.Range("M103").Value = -6833.5
$ws.Range("H132").Value = 3549.75
$ws.Range("I132").Value = 2400
$ws.Range("J132").Value = 6999
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 20997
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -26057
$ws.Range("H141").Value = 587299.6
$ws.Range("J141").Value = 587299.6
$ws.Range("L141").Value = 587299.6
$ws.Range("N141").Value = -597659.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 148.6
$ws.Range("J12").Value = 188.63637
$ws.Range("L12").Value = 565.9091100000001
$ws.Range("N12").Value = -911.9091100000001
$ws.Range("H40").Value = 121.35
$ws.Range("J40").Value = 188.25
$ws.Range("L40").Value = 753
$ws.Range("N40").Value = -891
$ws.Range("H122").Value = 919400.8
$ws.Range("J122").Value = 5052004.5
$ws.Range("L122").Value = 45468040.5
$ws.Range("N122").Value = -45472940.5
$ws.Range("H137").Value = 8302.9
$ws.Range("J137").Value = 10385.143
$ws.Range("L137").Value = 31155.429
$ws.Range("N137").Value = -41355.429
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18949.5
$ws.Range("I5").Value = 18949.5
$ws.Range("K5").Value = 18949.5
$ws.Range("M5").Value = -18837.5
$ws.Range("H70").Value = 43841.48
$ws.Range("I70").Value = 68292.94
$ws.Range("K70").Value = 68292.94
$ws.Range("M70").Value = -68022.94
$ws.Range("H73").Value = 43841.48
$ws.Range("I73").Value = 68292.94
$ws.Range("K73").Value = 68292.94
$ws.Range("M73").Value = -67356.94
$ws.Range("H122").Value = 14582.692
$ws.Range("J122").Value = 9889
$ws.Range("L122").Value = 29667
$ws.Range("N122").Value = -34567
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 930.27026
$ws.Range("I55").Value = 511.20834
$ws.Range("K55").Value = 511.20834
$ws.Range("M55").Value = -338.20834
$ws.Range("H82").Value = 910.3570999999999
$ws.Range("I82").Value = 685.4286
$ws.Range("J82").Value = 1135.2858
$ws.Range("K82").Value = 685.4286
$ws.Range("L82").Value = 1135.2858
$ws.Range("M82").Value = -324.4286
$ws.Range("N82").Value = -1857.2858
$ws.Range("H85").Value = 910.3570999999999
$ws.Range("I85").Value = 685.4286
$ws.Range("J85").Value = 1135.2858
$ws.Range("K85").Value = 685.4286
$ws.Range("L85").Value = 1135.2858
$ws.Range("M85").Value = 562.5714
$ws.Range("N85").Value = -3631.2858
$ws.Range("H93").Value = 1104.2
$ws.Range("I93").Value = 1011.5263
$ws.Range("J93").Value = 1264.2727
$ws.Range("K93").Value = 1011.5263
$ws.Range("L93").Value = 1264.2727
$ws.Range("M93").Value = 236.4737
$ws.Range("N93").Value = -3760.2727
$ws.Range("H96").Value = 149521
$ws.Range("J96").Value = 149521
$ws.Range("L96").Value = 149521
$ws.Range("N96").Value = -155013
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 7686333.5
$ws.Range("J5").Value = 7686333.5
$ws.Range("L5").Value = 7686333.5
$ws.Range("N5").Value = -7686557.5
$ws.Range("H15").Value = 72655
$ws.Range("J15").Value = 72655
$ws.Range("L15").Value = 72655
$ws.Range("N15").Value = -73231
$ws.Range("H107").Value = 1238.2285
$ws.Range("I107").Value = 808.875
$ws.Range("K107").Value = 2426.625
$ws.Range("M107").Value = -506.625
$ws.Range("H113").Value = 3333564.8
$ws.Range("I113").Value = 3333564.8
$ws.Range("K113").Value = 10000694.4
$ws.Range("M113").Value = -9998524.399999999
$ws.Range("H126").Value = 9518.191999999999
$ws.Range("I126").Value = 1896.4546
$ws.Range("J126").Value = 51437.75
$ws.Range("K126").Value = 5689.3638
$ws.Range("L126").Value = 154313.25
$ws.Range("M126").Value = -3219.3638
$ws.Range("N126").Value = -159253.25
$ws.Range("H132").Value = 1304.5883
$ws.Range("I132").Value = 1106.1538
$ws.Range("J132").Value = 1949.5
$ws.Range("K132").Value = 3318.4614
$ws.Range("L132").Value = 5848.5
$ws.Range("M132").Value = -788.4614000000001
$ws.Range("N132").Value = -10908.5

Write-Output "Applied all profit-sheet updates"